$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 2
$ws.Range("H2").Value = 247
$ws.Range("I2").Value = 166.66667
$ws.Range("J2").Value = 307.25
$ws.Range("K2").Value = 166.66667
$ws.Range("L2").Value = 307.25
$ws.Range("M2").Value = -53.66667000000001
$ws.Range("N2").Value = -533.25

# Row 43
$ws.Range("H43").Value = 3866.5
$ws.Range("I43").Value = 3100
$ws.Range("J43").Value = 4249.75
$ws.Range("K43").Value = 3100
$ws.Range("L43").Value = 4249.75
$ws.Range("M43").Value = -3031
$ws.Range("N43").Value = -4387.75

# Row 106
$ws.Range("H106").Value = 3878.8
$ws.Range("I106").Value = 4327
$ws.Range("K106").Value = 4327
$ws.Range("M106").Value = -3696

# Row 137
$ws.Range("H137").Value = 21327.637
$ws.Range("I137").Value = 1514.7142
$ws.Range("K137").Value = 4544.142599999999
$ws.Range("M137").Value = -1994.142599999999


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Range("H61").Value = 1986138.8
$ws.Range("I61").Value = 5059.4287
$ws.Range("K61").Value = 5059.4287
$ws.Range("M61").Value = -4847.4287

# Row 74
$ws.Range("H74").Value = 27667.578
$ws.Range("I74").Value = 1692.3
$ws.Range("K74").Value = 1692.3
$ws.Range("M74").Value = -818.3

# Row 77
$ws.Range("H77").Value = 27667.578
$ws.Range("I77").Value = 1692.3
$ws.Range("K77").Value = 8461.5
$ws.Range("M77").Value = -4093.5

# Row 95
$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492

# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("N96").ClearContents()

# Row 132
$ws.Range("H132").Value = 4140456.5
$ws.Range("I132").Value = 1222.963
$ws.Range("K132").Value = 3668.889
$ws.Range("M132").Value = -1138.889

# Row 136
$ws.Range("H136").Value = 1986138.8
$ws.Range("I136").Value = 5059.4287
$ws.Range("K136").Value = 15178.2861
$ws.Range("M136").Value = -12628.2861


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 34500.727
$ws.Range("I20").Value = 8777.556
$ws.Range("J20").Value = 52309.08
$ws.Range("K20").Value = 8777.556
$ws.Range("L20").Value = 52309.08
$ws.Range("M20").Value = -8530.556
$ws.Range("N20").Value = -52803.08

# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("N59").ClearContents()

# Row 97
$ws.Range("H97").Value = 7428
$ws.Range("I97").Value = 7428
$ws.Range("K97").Value = 7428
$ws.Range("M97").Value = -6437

# Row 134
$ws.Range("H134").Value = 54483.895
$ws.Range("I134").Value = 91194
$ws.Range("J134").Value = 22668.467
$ws.Range("K134").Value = 273582
$ws.Range("L134").Value = 68005.401
$ws.Range("M134").Value = -271047
$ws.Range("N134").Value = -73075.401


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 34743.5
$ws.Range("J31").Value = 41412.2
$ws.Range("L31").Value = 41412.2
$ws.Range("N31").Value = -42002.2

# Row 34
$ws.Range("H34").Value = 34743.5
$ws.Range("J34").Value = 41412.2
$ws.Range("L34").Value = 41412.2
$ws.Range("N34").Value = -41816.2

# Row 58
$ws.Range("H58").Value = 17762.75
$ws.Range("I58").Value = 8756.467000000001
$ws.Range("J58").Value = 28154.615
$ws.Range("K58").Value = 8756.467000000001
$ws.Range("L58").Value = 28154.615
$ws.Range("M58").Value = -8553.467000000001
$ws.Range("N58").Value = -28560.615

# Row 62
$ws.Range("H62").Value = 3299.5
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65
$ws.Range("H65").Value = 3299.5
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 132
$ws.Range("H132").Value = 44759830
$ws.Range("I132").Value = 2657.1538
$ws.Range("K132").Value = 7971.4614
$ws.Range("M132").Value = -5441.4614

# Row 134
$ws.Range("H134").Value = 32264836
$ws.Range("I134").Value = 2754.7646
$ws.Range("J134").Value = 71440216
$ws.Range("K134").Value = 8264.293799999999
$ws.Range("L134").Value = 214320648
$ws.Range("M134").Value = -5729.293799999999
$ws.Range("N134").Value = -214325718

# Row 136
$ws.Range("H136").Value = 17762.75
$ws.Range("I136").Value = 8756.467000000001
$ws.Range("J136").Value = 28154.615
$ws.Range("K136").Value = 26269.401
$ws.Range("L136").Value = 84463.845
$ws.Range("M136").Value = -23719.401
$ws.Range("N136").Value = -89563.845


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 56
$ws.Range("H56").Value = 7999
$ws.Range("I56").Value = 7999
$ws.Range("K56").Value = 7999
$ws.Range("M56").Value = -7469

# Row 60
$ws.Range("H60").Value = 3063.2
$ws.Range("I60").Value = 3420.6667
$ws.Range("J60").Value = 2527
$ws.Range("K60").Value = 10262.0001
$ws.Range("L60").Value = 7581
$ws.Range("M60").Value = -10011.0001
$ws.Range("N60").Value = -8083

# Row 122
$ws.Range("H122").Value = 10332.272
$ws.Range("J122").Value = 12560.111
$ws.Range("L122").Value = 113040.999
$ws.Range("N122").Value = -117940.999

# Row 141
$ws.Range("H141").Value = 6834
$ws.Range("I141").Value = 5499.25
$ws.Range("K141").Value = 16497.75
$ws.Range("M141").Value = -11317.75


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 4727.68
$ws.Range("I102").Value = 6189.625
$ws.Range("K102").Value = 6189.625
$ws.Range("M102").Value = -4567.625

# Row 122
$ws.Range("H122").Value = 1766.7693
$ws.Range("I122").Value = 1845.8948
$ws.Range("J122").Value = 1552
$ws.Range("K122").Value = 5537.6844
$ws.Range("L122").Value = 4656
$ws.Range("M122").Value = -3087.6844
$ws.Range("N122").Value = -9556

# Row 132
$ws.Range("H132").Value = 578033.6
$ws.Range("I132").Value = 5860.7
$ws.Range("K132").Value = 17582.1
$ws.Range("M132").Value = -15052.1


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 68
$ws.Range("H68").Value = 6972.2173
$ws.Range("I68").Value = 6064.1816
$ws.Range("J68").Value = 7804.5835
$ws.Range("K68").Value = 6064.1816
$ws.Range("L68").Value = 7804.5835
$ws.Range("M68").Value = -5315.1816
$ws.Range("N68").Value = -9302.583500000001

# Row 71
$ws.Range("H71").Value = 6972.2173
$ws.Range("I71").Value = 6064.1816
$ws.Range("J71").Value = 7804.5835
$ws.Range("K71").Value = 30320.908
$ws.Range("L71").Value = 39022.9175
$ws.Range("M71").Value = -26576.908
$ws.Range("N71").Value = -46510.9175

# Row 100
$ws.Range("H100").Value = 2407.5715
$ws.Range("I100").Value = 2283.5715
$ws.Range("J100").Value = 2655.5715
$ws.Range("K100").Value = 2283.5715
$ws.Range("L100").Value = 2655.5715
$ws.Range("M100").Value = -1742.5715
$ws.Range("N100").Value = -3737.5715

# Row 122
$ws.Range("H122").Value = 5280.9033
$ws.Range("J122").Value = 6133.5557
$ws.Range("L122").Value = 18400.6671
$ws.Range("N122").Value = -23300.6671

# Row 132
$ws.Range("H132").Value = 2589554.5
$ws.Range("I132").Value = 4280.55
$ws.Range("K132").Value = 12841.65
$ws.Range("M132").Value = -10311.65

# Row 136
$ws.Range("I136").Value = 18249.23
$ws.Range("K136").Value = 54747.69
$ws.Range("M136").Value = -52197.69


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 62
$ws.Range("H62").Value = 31499.666
$ws.Range("I62").Value = 31999.334
$ws.Range("J62").Value = 31000
$ws.Range("K62").Value = 31999.334
$ws.Range("L62").Value = 31000
$ws.Range("M62").Value = -31375.334
$ws.Range("N62").Value = -32248

# Row 65
$ws.Range("H65").Value = 31499.666
$ws.Range("I65").Value = 31999.334
$ws.Range("J65").Value = 31000
$ws.Range("K65").Value = 159996.67
$ws.Range("L65").Value = 155000
$ws.Range("M65").Value = -156876.67
$ws.Range("N65").Value = -161240

# Row 132
$ws.Range("H132").Value = 579765.5
$ws.Range("I132").Value = 2729.9375
$ws.Range("K132").Value = 8189.8125
$ws.Range("M132").Value = -5659.8125

# Row 136
$ws.Range("H136").Value = 702231.6
$ws.Range("I136").Value = 3166.5557
$ws.Range("K136").Value = 9499.667099999999
$ws.Range("M136").Value = -6949.667099999999

